# Add two new announcement rows to the table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by two rows first; this keeps the table's ref/autoFilter
# (and the underlying sheet dimension) in sync automatically.
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.ListRows.Add()
[void]$tbl.ListRows.Add()

# Copy the formatting (fill/border/number format/wrap) of the last
# pre-existing data row (row 3) down onto the two new rows (4 and 5).
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)

# Row 4: song practice reminder
$ws.Range("A4").Value = 45928
$ws.Range("B4").Value = "Song practice on 9/30/25 from 6:30pm to 8:30pm"
$ws.Range("C4").Value = "Muaj kawm nkauj hnub 9/30/25 thaum 6 moo 30 txog 8 moo 30"

# Row 5: dress-code reminder
$ws.Range("A5").Value = 45928
$ws.Range("B5").Value = "We will be wearing any shade of blue for the 40 years anniversary"
$ws.Range("C5").Value = "Peb hnav xim xiav rau lub 40 xyoo anniversary"

# Match the row heights (30pt, matching the wrapped-text rows above).
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30

# Move the active selection to C6, just like after typing the last row.
[void]$ws.Range("C6").Select()
